$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44435
$ws.Cells.Item(2, 13).Value = 130
$ws.Cells.Item(2, 14).Value = 1300
$ws.Cells.Item(2, 15).Value = 1300
$ws.Cells.Item(2, 16).Value = 1300
$ws.Cells.Item(2, 19).Value = 1300

$ws.Cells.Item(3, 4).Value = 45044
$ws.Cells.Item(3, 13).Value = 150
$ws.Cells.Item(3, 14).Value = 3500
$ws.Cells.Item(3, 15).Value = 3500
$ws.Cells.Item(3, 16).Value = 3500
$ws.Cells.Item(3, 19).Value = 3500

$ws.Cells.Item(4, 4).Value = 44432
$ws.Cells.Item(4, 13).Value = 30
$ws.Cells.Item(4, 14).Value = 1300
$ws.Cells.Item(4, 15).Value = 1300
$ws.Cells.Item(4, 16).Value = 1300
$ws.Cells.Item(4, 19).Value = 1300

$ws.Cells.Item(5, 4).Value = 44438
$ws.Cells.Item(5, 13).Value = 60
$ws.Cells.Item(5, 14).Value = 1200
$ws.Cells.Item(5, 15).Value = 1200
$ws.Cells.Item(5, 16).Value = 1200
$ws.Cells.Item(5, 19).Value = 1200

$ws.Cells.Item(6, 4).Value = 45041
$ws.Cells.Item(6, 13).Value = 80
$ws.Cells.Item(6, 14).Value = 3500
$ws.Cells.Item(6, 15).Value = 3500
$ws.Cells.Item(6, 16).Value = 3500
$ws.Cells.Item(6, 19).Value = 3500

$ws.Cells.Item(7, 4).Value = 44431
$ws.Cells.Item(7, 13).Value = 100
$ws.Cells.Item(7, 14).Value = 1300
$ws.Cells.Item(7, 15).Value = 1300
$ws.Cells.Item(7, 16).Value = 1300
$ws.Cells.Item(7, 19).Value = 1300

$ws.Cells.Item(8, 4).Value = 44418
$ws.Cells.Item(8, 13).Value = 40
$ws.Cells.Item(8, 14).Value = 1200
$ws.Cells.Item(8, 15).Value = 1200
$ws.Cells.Item(8, 16).Value = 1200
$ws.Cells.Item(8, 19).Value = 1200

$ws.Cells.Item(9, 4).Value = 45075
$ws.Cells.Item(9, 13).Value = 240
$ws.Cells.Item(9, 14).Value = 3200
$ws.Cells.Item(9, 15).Value = 3200
$ws.Cells.Item(9, 16).Value = 3200
$ws.Cells.Item(9, 19).Value = 3200

$ws.Cells.Item(10, 4).Value = 44748
$ws.Cells.Item(10, 13).Value = 300
$ws.Cells.Item(10, 14).Value = 2300
$ws.Cells.Item(10, 15).Value = 2300
$ws.Cells.Item(10, 16).Value = 2300
$ws.Cells.Item(10, 19).Value = 2300

$ws.Cells.Item(11, 4).Value = 44473
$ws.Cells.Item(11, 13).Value = 120
$ws.Cells.Item(11, 14).Value = 1200
$ws.Cells.Item(11, 15).Value = 1200
$ws.Cells.Item(11, 16).Value = 1200
$ws.Cells.Item(11, 19).Value = 1200

$ws.Cells.Item(12, 4).Value = 45068
$ws.Cells.Item(12, 13).Value = 50
$ws.Cells.Item(12, 14).Value = 3250
$ws.Cells.Item(12, 15).Value = 3250
$ws.Cells.Item(12, 16).Value = 3250
$ws.Cells.Item(12, 19).Value = 3250

$ws.Cells.Item(13, 4).Value = 44417
$ws.Cells.Item(13, 13).Value = 80
$ws.Cells.Item(13, 14).Value = 1200
$ws.Cells.Item(13, 15).Value = 1200
$ws.Cells.Item(13, 16).Value = 1200
$ws.Cells.Item(13, 19).Value = 1200

$ws.Cells.Item(14, 4).Value = 44830
$ws.Cells.Item(14, 13).Value = 50
$ws.Cells.Item(14, 14).Value = 2500
$ws.Cells.Item(14, 15).Value = 2500
$ws.Cells.Item(14, 16).Value = 2500
$ws.Cells.Item(14, 19).Value = 2500

$ws.Cells.Item(15, 4).Value = 44405
$ws.Cells.Item(15, 13).Value = 50
$ws.Cells.Item(15, 14).Value = 1200
$ws.Cells.Item(15, 15).Value = 1200
$ws.Cells.Item(15, 16).Value = 1200
$ws.Cells.Item(15, 19).Value = 1200

$ws.Cells.Item(16, 4).Value = 44476
$ws.Cells.Item(16, 13).Value = 80
$ws.Cells.Item(16, 14).Value = 1200
$ws.Cells.Item(16, 15).Value = 1200
$ws.Cells.Item(16, 16).Value = 1200
$ws.Cells.Item(16, 19).Value = 1200

$ws.Cells.Item(17, 4).Value = 44357
$ws.Cells.Item(17, 13).Value = 35
$ws.Cells.Item(17, 14).Value = 1000
$ws.Cells.Item(17, 15).Value = 1000
$ws.Cells.Item(17, 16).Value = 1000
$ws.Cells.Item(17, 19).Value = 1000

$ws.Cells.Item(18, 4).Value = 44811
$ws.Cells.Item(18, 13).Value = 60
$ws.Cells.Item(18, 14).Value = 2500
$ws.Cells.Item(18, 15).Value = 2500
$ws.Cells.Item(18, 16).Value = 2500
$ws.Cells.Item(18, 19).Value = 2500

$ws.Cells.Item(19, 4).Value = 44763
$ws.Cells.Item(19, 13).Value = 50
$ws.Cells.Item(19, 14).Value = 2300
$ws.Cells.Item(19, 15).Value = 2300
$ws.Cells.Item(19, 16).Value = 2300
$ws.Cells.Item(19, 19).Value = 2300

$ws.Cells.Item(20, 4).Value = 44760
$ws.Cells.Item(20, 13).Value = 80
$ws.Cells.Item(20, 14).Value = 2300
$ws.Cells.Item(20, 15).Value = 2300
$ws.Cells.Item(20, 16).Value = 2300
$ws.Cells.Item(20, 19).Value = 2300

$ws.Cells.Item(21, 4).Value = 44343
$ws.Cells.Item(21, 13).Value = 60
$ws.Cells.Item(21, 14).Value = 1300
$ws.Cells.Item(21, 15).Value = 1300
$ws.Cells.Item(21, 16).Value = 1300
$ws.Cells.Item(21, 19).Value = 1300

$ws.Cells.Item(22, 4).Value = 44424
$ws.Cells.Item(22, 13).Value = 50
$ws.Cells.Item(22, 14).Value = 1200
$ws.Cells.Item(22, 15).Value = 1200
$ws.Cells.Item(22, 16).Value = 1200
$ws.Cells.Item(22, 19).Value = 1200

$ws.Cells.Item(23, 4).Value = 45055
$ws.Cells.Item(23, 13).Value = 25
$ws.Cells.Item(23, 14).Value = 2800
$ws.Cells.Item(23, 15).Value = 2800
$ws.Cells.Item(23, 16).Value = 2800
$ws.Cells.Item(23, 19).Value = 2800

$ws.Cells.Item(24, 4).Value = 45042
$ws.Cells.Item(24, 13).Value = 25
$ws.Cells.Item(24, 14).Value = 3500
$ws.Cells.Item(24, 15).Value = 3500
$ws.Cells.Item(24, 16).Value = 3500
$ws.Cells.Item(24, 19).Value = 3500

$ws.Cells.Item(25, 4).Value = 45062
$ws.Cells.Item(25, 13).Value = 60
$ws.Cells.Item(25, 14).Value = 3200
$ws.Cells.Item(25, 15).Value = 3200
$ws.Cells.Item(25, 16).Value = 3200
$ws.Cells.Item(25, 19).Value = 3200

$ws.Cells.Item(26, 4).Value = 44749
$ws.Cells.Item(26, 13).Value = 120
$ws.Cells.Item(26, 14).Value = 2300
$ws.Cells.Item(26, 15).Value = 2300
$ws.Cells.Item(26, 16).Value = 2300
$ws.Cells.Item(26, 19).Value = 2300

$ws.Cells.Item(27, 4).Value = 45054
$ws.Cells.Item(27, 13).Value = 25
$ws.Cells.Item(27, 14).Value = 2500
$ws.Cells.Item(27, 15).Value = 2500
$ws.Cells.Item(27, 16).Value = 2500
$ws.Cells.Item(27, 19).Value = 2500

$ws.Cells.Item(28, 4).Value = 44762
$ws.Cells.Item(28, 13).Value = 50
$ws.Cells.Item(28, 14).Value = 2300
$ws.Cells.Item(28, 15).Value = 2300
$ws.Cells.Item(28, 16).Value = 2300
$ws.Cells.Item(28, 19).Value = 2300

$ws.Cells.Item(29, 4).Value = 45076
$ws.Cells.Item(29, 13).Value = 100
$ws.Cells.Item(29, 14).Value = 2600
$ws.Cells.Item(29, 15).Value = 2600
$ws.Cells.Item(29, 16).Value = 2600
$ws.Cells.Item(29, 19).Value = 2600

$ws.Cells.Item(30, 4).Value = 44753
$ws.Cells.Item(30, 13).Value = 160
$ws.Cells.Item(30, 14).Value = 2300
$ws.Cells.Item(30, 15).Value = 2300
$ws.Cells.Item(30, 16).Value = 2300
$ws.Cells.Item(30, 19).Value = 2300

$ws.Cells.Item(31, 4).Value = 44812
$ws.Cells.Item(31, 13).Value = 50
$ws.Cells.Item(31, 14).Value = 2500
$ws.Cells.Item(31, 15).Value = 2500
$ws.Cells.Item(31, 16).Value = 2500
$ws.Cells.Item(31, 19).Value = 2500
